$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "ssurgwsoadev4-oci.opc.oracleoutsourcing.com"
$ws.Range("B3").Value = "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do"
$ws.Range("C3").Value = "su"
$ws.Range("D3").Value = "gw"

$ws.Hyperlinks.Add($ws.Range("B3"), "https://ssurgwsoadev4-oci.opc.oracleoutsourcing.com/pc/PolicyCenter.do") | Out-Null
$ws.Range("B3").Style = $ws.Range("B2").Style

$ws.Range("B8").Select() | Out-Null
